# ajout dans tableau de bord
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# New row 19 mirrors the formatting of row 18 (same row height category, date style,
# name style, and wrapped-text styles for the functionality / contribution columns).
$ws.Range("B18").Copy()
$ws.Range("B19").PasteSpecial($xlPasteFormats)

$ws.Range("C18").Copy()
$ws.Range("C19").PasteSpecial($xlPasteFormats)

$ws.Range("D18").Copy()
$ws.Range("D19").PasteSpecial($xlPasteFormats)

$ws.Range("E18").Copy()
$ws.Range("E19").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# Fill in the new contribution entry
$ws.Range("B19").Value = (Get-Date -Year 2021 -Month 1 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("C19").Value = "Thomas"
$ws.Range("D19").Value = "11`n12`n13"
$ws.Range("E19").Value = "Tout`nTout`nTout"

$ws.Rows.Item(19).RowHeight = 45

$ws.Range("E20").Select()

